$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# --- Refresh the "time_taken" timestamps on the existing "data" sheet ---
$dataWs.Range("F2").Value = "2021-10-05 14:20:54.631097"
$dataWs.Range("F3").Value = "2021-10-05 14:20:54.631105"
$dataWs.Range("F4").Value = "2021-10-05 14:20:54.631108"
$dataWs.Range("F5").Value = "2021-10-05 14:20:54.631111"
$dataWs.Range("F6").Value = "2021-10-05 14:20:54.631114"
$dataWs.Range("F7").Value = "2021-10-05 14:20:54.631117"
$dataWs.Range("F8").Value = "2021-10-05 14:20:54.631119"
$dataWs.Range("F9").Value = "2021-10-05 14:20:54.631122"
$dataWs.Range("F10").Value = "2021-10-05 14:20:54.631125"
$dataWs.Range("F11").Value = "2021-10-05 14:20:54.631127"
$dataWs.Range("F12").Value = "2021-10-05 14:20:54.631130"
$dataWs.Range("F13").Value = "2021-10-05 14:20:54.631133"
$dataWs.Range("F14").Value = "2021-10-05 14:20:54.631135"
$dataWs.Range("F15").Value = "2021-10-05 14:20:54.631138"
$dataWs.Range("F16").Value = "2021-10-05 14:20:54.631140"
$dataWs.Range("F17").Value = "2021-10-05 14:20:54.631143"
$dataWs.Range("F18").Value = "2021-10-05 14:20:54.631146"
$dataWs.Range("F19").Value = "2021-10-05 14:20:54.631149"
$dataWs.Range("F20").Value = "2021-10-05 14:20:54.631151"
$dataWs.Range("F21").Value = "2021-10-05 14:20:54.631154"
$dataWs.Range("F22").Value = "2021-10-05 14:20:54.631157"
$dataWs.Range("F23").Value = "2021-10-05 14:20:54.631159"
$dataWs.Range("F24").Value = "2021-10-05 14:20:54.631162"
$dataWs.Range("F25").Value = "2021-10-05 14:20:54.631165"
$dataWs.Range("F26").Value = "2021-10-05 14:20:54.631168"
$dataWs.Range("F27").Value = "2021-10-05 14:20:54.631171"
$dataWs.Range("F28").Value = "2021-10-05 14:20:54.631173"
$dataWs.Range("F29").Value = "2021-10-05 14:20:54.631176"
$dataWs.Range("F30").Value = "2021-10-05 14:20:54.631179"
$dataWs.Range("F31").Value = "2021-10-05 14:20:54.631181"
$dataWs.Range("F32").Value = "2021-10-05 14:20:54.631184"
$dataWs.Range("F33").Value = "2021-10-05 14:20:54.631187"
$dataWs.Range("F34").Value = "2021-10-05 14:20:54.631190"
$dataWs.Range("F35").Value = "2021-10-05 14:20:54.631192"
$dataWs.Range("F36").Value = "2021-10-05 14:20:54.631195"
$dataWs.Range("F37").Value = "2021-10-05 14:20:54.631198"
$dataWs.Range("F38").Value = "2021-10-05 14:20:54.631200"
$dataWs.Range("F39").Value = "2021-10-05 14:20:54.631203"
$dataWs.Range("F40").Value = "2021-10-05 14:20:54.631206"
$dataWs.Range("F41").Value = "2021-10-05 14:20:54.631208"
$dataWs.Range("F42").Value = "2021-10-05 14:20:54.631212"
$dataWs.Range("F43").Value = "2021-10-05 14:20:54.631214"
$dataWs.Range("F44").Value = "2021-10-05 14:20:54.631217"
$dataWs.Range("F45").Value = "2021-10-05 14:20:54.631219"
$dataWs.Range("F46").Value = "2021-10-05 14:20:54.631222"
$dataWs.Range("F47").Value = "2021-10-05 14:20:54.631225"
$dataWs.Range("F48").Value = "2021-10-05 14:20:54.631228"

# --- Add the new "metadata" sheet right after "data" ---
$ws = $wb.Worksheets.Add($null, $dataWs)
$ws.Name = "metadata"

# Copy the header formatting (bold font + thin border + center/top alignment)
# from the "data" sheet's header row so the new header cells reuse the same
# style instead of creating a brand-new one.
$dataWs.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$dataWs.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Hypogonadotropic hypogonadism"
$ws.Range("C2").Value = 92
$ws.Range("D2").Value = "'1.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "2021-07-19T08:58:37.489830Z"
$ws.Range("F2").Value = "2021-10-05 14:20:54.627847"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/92/?format=json"
